$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "NOT SUPPORTED YET"
$ws.Range("D4").Select()
